$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct accented names used as simulation/column headers
$ws.Range("K1").Value = "Angela"
$ws.Range("M1").Value = "Rocio"

# Remove accents from adjective pair labels (rows used by the "wimpgrid" simulation titles)
$ws.Range("O13").Value = "Energico"
$ws.Range("A14").Value = "Egocentrico"
$ws.Range("O14").Value = "Empatico"
$ws.Range("A17").Value = "Frio"
$ws.Range("O17").Value = "Calido"
$ws.Range("A19").Value = "Antipatico"
$ws.Range("O20").Value = "Grunyon"
